# Update Leave Card File 1/12/2024 4:30 pm
#
# On the "2018 LEAVE CREDITS" sheet, Table15 gets a new undertime row
# inserted above the old row 75 (pushing the "2023"/"2024" section
# markers and every later row down by one, and extending the table from
# A8:K135 to A8:K136). Two new undertime entries are also recorded:
#   - row 71 (date 44805): "UT(0-6-0)" for 0.75 day
#   - new row 75:          "UT(3-6-0)" for 3.75 days
# The CONVERTION sheet's daily-earn calculator inputs are updated from
# 5 hours / 55 minutes to 6 hours / 0 minutes.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("2018 LEAVE CREDITS")

# Insert a new table row above (old) row 75 - this shifts old rows
# 75..135 down to 76..136, matching the way Excel's "Insert Table Rows
# Above" pushes the rest of the table (and the trailing total-style row)
# down by one.
$ws.Rows("75:75").Insert()

# The inserted row starts out with generic formatting; copy the format
# of the row above (old row 74, a normal data row) onto it so its style
# matches the rest of the table.
$ws.Range("A74:K74").Copy()
$ws.Range("A75:K75").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Grow the Table15 listobject to include the freshly inserted row.
$lo = $ws.ListObjects.Item("Table15")
$lo.Resize($ws.Range("A8:K136"))

# Populate the new row: PARTICULARS = "UT(3-6-0)", Absence Undertime W/
# Pay = 3.75, and the calculated "EARNED " helper column formula.
$ws.Range("B75").Value2 = "UT(3-6-0)"
$ws.Range("D75").Value2 = 3.75
$ws.Range("G75").Formula = '=IF(ISBLANK(Table15[[#This Row],[EARNED]]),"",Table15[[#This Row],[EARNED]])'

# Make sure the final (totals-style) row keeps the same calculated
# column formula shape as the rest of the table.
$ws.Range("G136").Formula = '=IF(ISBLANK(Table15[[#This Row],[EARNED]]),"",Table15[[#This Row],[EARNED]])'

# Record a new undertime entry on the existing (still) row 71.
$ws.Range("B71").Value2 = "UT(0-6-0)"
$ws.Range("D71").Value2 = 0.75

# Update the CONVERTION sheet's "daily earn calculator" hour/minute
# inputs: 5h55m -> 6h0m.
$ws2 = $wb.Worksheets.Item("CONVERTION")
$ws2.Range("E3").Value2 = 6
$ws2.Range("F3").ClearContents()

# Restore the sheet's active selection to match where the edit left off.
$ws.Activate()
$ws.Range("F79").Select()
